$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settlements")

# Match the formatting already used by column B (style index 1)
$ws.Range("B1:B4").Copy()
$ws.Range("C1:C4").PasteSpecial(-4122)

$ws.Range("C1").Value = "1d2"
$ws.Range("C2").Value = "1d4+1"
$ws.Range("C3").Value = "1d6+2"
$ws.Range("C4").Value = "1d8+3"
